$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date (Excel serial 45192 = 2023-09-23) on every
# data row (2-230). Update it to serial 45202 (2023-10-03).
$ws.Range("C2:C230").Value = 45202
